$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in @(1,2,3,4,6,7)) {
    $src = $ws.Range("C$r")
    $dst = $ws.Range("D$r")
    $src.Copy()
    $dst.PasteSpecial(-4122)
    $dst.Value = "fatta"
}
$excel.CutCopyMode = $false

$ws.Range("D8").Font.Underline = $true

$ws.Range("D8").Select()

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
